# Apply the commit's changes to the workbook:
#  - Sheet1 header row: PropAddr/PropCity/PropState/BPO/PropZip/Comment -> A/B/C/D/E/F
#  - Analysis header row: State/BPO/Calc -> A/B/Calc (keep Calc last)
#  - Analysis formulas now reference Sheet1!A.. / Sheet1!B.. (instead of C/D) and
#    C = A*B (instead of 2*B)
#  - Active sheet becomes "Analysis" (tabSelected moves from Sheet1 to Analysis)
#  - Selection: Sheet1 -> F1, Analysis -> D8

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Analysis")

# --- Sheet1 header row values ---
$ws1.Range("A1").Value = "A"
$ws1.Range("B1").Value = "B"
$ws1.Range("C1").Value = "C"
$ws1.Range("D1").Value = "D"
$ws1.Range("E1").Value = "E"
$ws1.Range("F1").Value = "F"

# --- Analysis header row values ---
$ws2.Range("A1").Value = "A"
$ws2.Range("B1").Value = "B"
$ws2.Range("C1").Value = "Calc"

# --- Analysis formulas ---
$ws2.Range("A2").Formula = "=Sheet1!A2"
$ws2.Range("B2").Formula = "=Sheet1!B2"
$ws2.Range("C2").Formula = "=A2*B2"

$ws2.Range("A3").Formula = "=Sheet1!A3"
$ws2.Range("B3").Formula = "=Sheet1!B3"
$ws2.Range("C3").Formula = "=A3*B3"

# --- Selections on each sheet ---
$ws1.Range("F1").Select() | Out-Null
$ws2.Range("D8").Select() | Out-Null

# --- Make Analysis the active/selected sheet (tab) ---
$ws2.Activate() | Out-Null
